# "added failsafe for email body on smtp"
#
# The sheet used to hold a big mailing list (TONAME / EMAIL columns with
# dozens of first names, last names and generated gmail addresses). The
# author trimmed it down to just the first two sample rows and pointed the
# EMAIL column at a single real, failsafe address instead - turning H2/H3
# into live mailto: hyperlinks so a misconfigured SMTP "to" can't silently
# fire off to the whole old test list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$failsafeEmail = "shahimrans64@gmail.com"

# Drop all the old sample rows (4-28) in the TONAME (C/G) and EMAIL (H)
# columns - ClearContents removes the now-empty rows/cells entirely and the
# sheet's dimension + shared-string table shrink to match.
$ws.Range("C4:C8").ClearContents()
$ws.Range("G4:G8").ClearContents()
$ws.Range("H4:H28").ClearContents()

# Point the two remaining EMAIL cells at the single failsafe address.
$ws.Range("H2").Value = $failsafeEmail
$ws.Range("H3").Value = $failsafeEmail

# Remove any pre-existing hyperlinks on those cells before re-adding, so
# re-running this script stays idempotent.
$null = $ws.Range("H2").Hyperlinks.Delete()
$null = $ws.Range("H3").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:" + $failsafeEmail)
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:" + $failsafeEmail)

# Match the author's last-saved selection.
$null = $ws.Range("H3").Select()

Write-Output "applied smtp failsafe email edit"
